# Generate Report for Handback
# - Update the "Ready for handoff" status (Overview sheet) to "Handback transform failed"
# - Populate the "Error Detail" column (P) on the zh-cn and de-de sheets with the
#   handback/handoff file-name mismatch message
# - Widen the "Error Detail" column (P) on both locale sheets to fit the new text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The shared "Ready for handoff" status text is used for the e6f39f98 row on the
# Overview sheet (both the Status and duplicated-status columns) as well as the
# Status column on each locale sheet, and must now reflect that the handback
# transform failed.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error details for the failed handback of the e6f39f98 file, per locale.
$zhcn.Range("P3").Value = "Handback file name: ub5q0ipg.naw is different with handoff file name: e6f39f98-32b8-417c-be47-5a1ec7d811f0.3fdfc5a585b18e25c10e31fae7b069eb291b1f89.zh-cn."
$dede.Range("P3").Value = "Handback file name: ub5q0ipg.naw is different with handoff file name: e6f39f98-32b8-417c-be47-5a1ec7d811f0.3fdfc5a585b18e25c10e31fae7b069eb291b1f89.de-de."

# Widen the Error Detail column (P, the 16th column) so the new message is readable.
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
